# Fix total marks error on the marksheet.
# Marking scheme changed: right-answer marks 5 -> 4, wrong-answer penalty -1 -> -2.
# This updates the computed totals and the "score / max" summary text accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": marks per right answer and penalty per wrong answer
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": recomputed totals using the corrected marking scheme
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "82 / 112"
